$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.59"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.93"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.247"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05857"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.460"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.335"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8085"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8787"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1377"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07253"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03071"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03052"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09323"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001544"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04710"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006001"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006221"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001264"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004574"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008700"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.177"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03783"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006350"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1054"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002560"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007804"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005494"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5401"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02174"
